# This workbook is a crypto-price tracker. The automation periodically refreshes
# the Price (column D) and Volume(1h) (column E) figures for each listed coin, and
# occasionally re-ranks rows (e.g. rows 40/41 below swap places as 'Cosmos' and
# 'dogwifhat' trade ranking positions). This run applies the latest scrape values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of cell address -> refreshed value taken from the latest data pull.
$updates = [ordered]@{
    'D2' = '62.694.83'
    'E2' = '  -1.51%  '
    'D3' = '3.026.17'
    'E3' = '  -1.80%  '
    'E4' = '  +0.00%  '
    'D5' = '583.80'
    'E5' = '  -1.07%  '
    'D6' = '148.80'
    'E6' = '  -4.49%  '
    'E7' = '  +0.00%  '
    'E8' = '  -3.30%  '
    'D9' = '3.023.37'
    'E9' = '  -1.85%  '
    'E10' = '  -3.21%  '
    'D11' = '5.68'
    'E11' = '  -2.66%  '
    'E12' = '  -2.12%  '
    'E13' = '  -3.50%  '
    'D14' = '35.37'
    'E14' = '  -4.91%  '
    'E15' = '  +1.82%  '
    'D16' = '3.534.79'
    'E16' = '  -1.62%  '
    'D17' = '7.07'
    'E17' = '  -1.27%  '
    'D18' = '62.729.87'
    'E18' = '  -1.43%  '
    'D19' = '3.029.05'
    'E19' = '  -1.70%  '
    'D20' = '468.53'
    'E20' = '  -1.90%  '
    'D21' = '14.04'
    'E21' = '  -3.02%  '
    'D22' = '0.692'
    'E22' = '  -2.55%  '
    'D23' = '7.42'
    'E23' = '  -1.61%  '
    'D24' = '2.37'
    'E24' = '  -1.55%  '
    'D25' = '80.97'
    'E25' = '  -0.14%  '
    'D26' = '12.42'
    'E26' = '  -3.17%  '
    'D27' = '10.43'
    'E27' = '  +0.91%  '
    'E28' = '  -0.06%  '
    'E29' = '  +0.07%  '
    'D30' = '7.26'
    'E30' = '  -3.58%  '
    'E31' = '  -1.36%  '
    'D32' = '2.14'
    'E32' = '  -1.30%  '
    'D33' = '27.44'
    'E33' = '  +0.87%  '
    'E34' = '  -4.60%  '
    'E35' = '  -1.49%  '
    'D36' = '0.0₃0796'
    'E36' = '  -5.56%  '
    'D37' = '5.77'
    'E37' = '  -4.44%  '
    'D38' = '2.15'
    'E38' = '  -2.48%  '
    'D39' = '50.24'
    'E39' = '  -1.09%  '
    'B40' = 'Cosmos'
    'C40' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D40' = '9.00'
    'E40' = '  -3.90%  '
    'B41' = 'dogwifhat'
    'C41' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D41' = '2.95'
    'E41' = '  -12.98%  '
    'D42' = '423.05'
    'E42' = '  -4.96%  '
    'E43' = '  -1.88%  '
    'E44' = '  +0.65%  '
    'D45' = '2.805.48'
    'E45' = '  +0.00%  '
    'E46' = '  -1.41%  '
    'D47' = '37.77'
    'E47' = '  -8.37%  '
    'D48' = '129.09'
    'E48' = '  -2.06%  '
    'D49' = '0.999'
    'E49' = '  -0.05%  '
    'D50' = '24.39'
    'E50' = '  -3.60%  '
    'E51' = '  -1.20%  '
}

foreach ($addr in $updates.Keys) {
    $range = $ws.Range($addr)
    $newValue = $updates[$addr]
    # Prefix with an apostrophe so Excel stores the refreshed figure as literal text
    # (matching the sheet's existing text-formatted Price/Volume/Coin/Link columns)
    # instead of silently re-interpreting strings such as '583.80' or '0.999' as
    # numbers and dropping their displayed trailing zeros.
    $range.Value = "'" + $newValue
    # Restore the default 'Normal' style so no incidental number formatting (from
    # the text-forcing prefix above) lingers on the cell.
    $range.Style = 'Normal'
}
